# Update column C (dates) from 45207 (2023-10-08) to 45208 (2023-10-09)
# for all data rows (2 through 54) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 54
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
